$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B26").Value = 0.9999674344715328
$ws.Range("C2:C26").Value = 0.9989400190852257
$ws.Range("D2:D26").Value = 0.9999999791127612
$ws.Range("E2:E26").Value = 0.9999790395827924
$ws.Range("F2:F26").Value = 0.9999952325641216
$ws.Range("G2:G26").Value = 0.0000303984913626197098
$ws.Range("H2:H26").Value = 0.0009894456561571146
$ws.Range("I2:I26").Value = 0.0000000161730318727365908
$ws.Range("J2:J26").Value = 0.00000547499558639807835
$ws.Range("K2:K26").Value = 0.00000274558430913540686
$ws.Range("L2:L26").Value = 0.000349036811621329
$ws.Range("M2:M26").Value = 0.005513482689065026
$ws.Range("N2:N26").Value = 1.000060120975632
$ws.Range("O2:O26").Value = 0.005748203061621094
$ws.Range("P2:P26").Value = 94.80223515402628
$ws.Range("Q2:Q26").Value = 139.9006406741497
